function Force-Split($doc, $pos) {
    # Forces a run boundary at a given character position by briefly adding
    # and then removing a bookmark there -- the engine always serializes a
    # bookmark as a standalone element, which splits any run it lands inside.
    $r = $doc.Range($pos, $pos)
    $doc.Bookmarks.Add("TempSplitMark", $r)
    $doc.Bookmarks("TempSplitMark").Delete()
}

$d = $word.ActiveDocument

# =====================================================================
# Change 1: "...with its loose typing system..." -> "...with its dynamic
# typing system..." as three separate runs (unchanged lead-in sentence,
# "dynamic" standalone, trailing clause).
# =====================================================================

$rngLoose = $d.Content
$rngLoose.Find.Execute("loose", $true, $false, $false, $false, $false, $true, 1, $false, "dynamic", 2)

$rngItsForemost = $d.Content
$rngItsForemost.Find.Execute("Its foremost")
Force-Split $d $rngItsForemost.Start

$rngDynamic = $d.Content
$rngDynamic.Find.Execute("dynamic")
Force-Split $d $rngDynamic.Start
Force-Split $d $rngDynamic.End

# =====================================================================
# Change 2: "...iterators, and map, filter..." -> "...iterators, and
# decorator, map, filter..." with the _GoBack bookmark relocated to sit
# right before "map" (after "decorator, ") instead of right after
# "create". Also the spell-check markup around "Indra" goes away since
# that part of the paragraph is being retyped/merged.
# =====================================================================

# Re-typing the "context of the Indra ABM" span (round-tripped through a
# placeholder) drops the spellStart/spellEnd proofErr wrapper around
# "Indra" and lets the surrounding text re-merge into one run.
$rngIndra = $d.Content
$rngIndra.Find.Execute("context of the Indra ABM")
$rngIndra.Text = "ZZZPLACEHOLDERZZZ"
$rngIndra2 = $d.Content
$rngIndra2.Find.Execute("ZZZPLACEHOLDERZZZ")
$rngIndra2.Text = "context of the Indra ABM"

# Insert "decorator, " right before "map, filter".
$rngMap = $d.Content
$rngMap.Find.Execute("map, filter")
$rngMap.Collapse(1)
$rngMap.InsertBefore("decorator, ")

# Move the (hidden) _GoBack bookmark from after "create" to right before
# "map" (i.e. right after the newly inserted "decorator, ").
$rngMap2 = $d.Content
$rngMap2.Find.Execute("map, filter")
$newBmRng = $d.Range($rngMap2.Start, $rngMap2.Start)
$d.Bookmarks.Add("_GoBack", $newBmRng)

# Force the run boundaries needed so "create" stays its own run, "This
# talk will present..." stays separate from "context of the Indra...",
# and "...and " stays separate from "decorator, ".
$rngCreate = $d.Content
$rngCreate.Find.Execute("create")
Force-Split $d $rngCreate.Start

$rngContext = $d.Content
$rngContext.Find.Execute("context of the")
Force-Split $d $rngContext.Start

$rngDecorator = $d.Content
$rngDecorator.Find.Execute("decorator, ")
Force-Split $d $rngDecorator.Start

Write-Output "edit complete"
